$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 88

# Column A holds a date-like string ("2025/10/10") that must stay plain text
# (matches the rest of the sheet, which stores dates as text, not real dates).
# Temporarily force text format so the value isn't auto-converted into a date
# serial, then clear the formatting again so the cell keeps the workbook's
# default (unstyled) look, same as every other data row.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/10"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "金"
$ws.Cells.Item($row, 3).Value = 9
$ws.Cells.Item($row, 4).Value = 30
